# Project Euler 004 A.xlsx — "feat: improve presentation details"
#
# Summary of the change:
#  - The whole content block (rows 3-21) shifts up by one row (rows 1-2 were
#    always empty, so deleting one of the leading empty rows accomplishes
#    this and lets every formula/row reference auto-adjust).
#  - A couple of text fixes: "Projet Euler" -> "Project Euler", the two
#    section headers get reworded, and a trailing "(end)" marker is added
#    on what used to be a blank spacer row at the bottom.
#  - The hyperlink that was anchored on the title-url row needs to be
#    re-anchored on its new row.
#  - Cosmetic: selection moves to B3, zoom normalizes to 100%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift every row up by one -----------------------------------------
# Rows 1 and 2 are blank in the original sheet (content starts at row 3), so
# removing row 1 slides rows 3..21 up to 2..20 and every same-row formula
# reference (B12 inside the C12 formula, etc.) is kept in sync by the engine.
$ws.Rows.Item(1).Delete()

# --- 2. Re-anchor the hyperlink on its new row (B5 -> B4) ------------------
$ws.Range("B5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B4"), "https://projecteuler.net/problem=4") | Out-Null

# --- 3. Text fixes ----------------------------------------------------------
$ws.Range("B2").Value = "Project Euler 4: Largest Palindrome Product"
$ws.Range("B9").Value = "1) One-liner based on recursion"
$ws.Range("B14").Value = "2) One-liner based on array formulas"

# --- 4. New "(end)" marker on the last spacer row --------------------------
$ws.Range("B19").Value = "(end)"
$ws.Range("B19").Font.ThemeFont = 0

# --- 5. Restore the single-cell dynamic-array formula markers --------------
# Shifting the rows makes the engine re-emit these four formulas as plain
# (non-array) formulas; re-applying them as array formulas restores the
# t="array" ref="..." marking Excel uses for LAMBDA/LET one-liners.
foreach ($addr in @("C11", "C12", "C16", "C17")) {
    $cell = $ws.Range($addr)
    $cell.FormulaArray = $cell.Formula
}

# --- 6. View cosmetics -------------------------------------------------------
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 100
